$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a brand-new column before column C. Excel shifts the old
# C:F ("Cc","e0","k" + the trailing numeric column) one column to the
# right, so the old column C/D/E/F become D/E/F/G.
$ws.Columns("C:C").Insert()

# Header for the freshly inserted column.
$ws.Range("C1").Value = "Surcharge"

# The values that used to live in column F (now shifted to G) are really
# the new "Surcharge" column's data, so move them into C2:C85 and give
# them the same number format already used by the other numeric columns.
$ws.Range("C2:C85").Value2 = $ws.Range("G2:G85").Value2
$ws.Range("C2:C85").NumberFormat = $ws.Range("D2:D85").NumberFormat

# Drop the now-duplicated trailing column.
$ws.Columns("G:G").Delete()

# Reproduce the saved selection state (columns C:F selected, active cell C1).
[void]$ws.Columns("C:F").Select()
